$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new SKU ("G713IC-HX008") was added to the product-codes list. Insert a
# fresh row above the current row 2 so every code below it shifts down one
# row, then fill the new row with the SKU.
$ws.Range("A2").EntireRow.Insert()
$ws.Range("A2").Value = "G713IC-HX008"

# Leave the selection where the edit happened, matching the saved view state.
$ws.Range("C11").Select()
